# Auto-generated edit script: updates cached market/profit values
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (columns H-N),
# matching refreshed currentAveragePrice / LevePrice / LeveProfit figures.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1358.125
$ws.Range("I41").Value = 1509.5714
$ws.Range("J41").Value = 298
$ws.Range("K41").Value = 1509.5714
$ws.Range("L41").Value = 298
$ws.Range("M41").Value = -1069.5714
$ws.Range("N41").Value = -1178
$ws.Range("H82").Value = 8769.833000000001
$ws.Range("I82").Value = 3872.3333
$ws.Range("J82").Value = 13667.333
$ws.Range("K82").Value = 11616.9999
$ws.Range("L82").Value = 41001.999
$ws.Range("M82").Value = -11210.9999
$ws.Range("N82").Value = -41813.999
$ws.Range("H85").Value = 8769.833000000001
$ws.Range("I85").Value = 3872.3333
$ws.Range("J85").Value = 13667.333
$ws.Range("K85").Value = 11616.9999
$ws.Range("L85").Value = 41001.999
$ws.Range("M85").Value = -10212.9999
$ws.Range("N85").Value = -43809.999
$ws.Range("H129").Value = 4420.92
$ws.Range("I129").Value = 1119.9
$ws.Range("J129").Value = 6621.6
$ws.Range("K129").Value = 3359.7
$ws.Range("L129").Value = 19864.8
$ws.Range("M129").Value = 1640.3
$ws.Range("N129").Value = -29864.8
$ws.Range("H137").Value = 4508.7104
$ws.Range("I137").Value = 3892.2273
$ws.Range("J137").Value = 5356.375
$ws.Range("K137").Value = 11676.6819
$ws.Range("L137").Value = 16069.125
$ws.Range("M137").Value = -9126.6819
$ws.Range("N137").Value = -21169.125
$ws.Range("H138").Value = 8176.8184
$ws.Range("I138").Value = 8502.833000000001
$ws.Range("J138").Value = 8104.3706
$ws.Range("K138").Value = 25508.499
$ws.Range("L138").Value = 24313.1118
$ws.Range("M138").Value = -20368.499
$ws.Range("N138").Value = -34593.1118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3303.5417
$ws.Range("I2").Value = 3116.8333
$ws.Range("J2").Value = 3863.6667
$ws.Range("K2").Value = 3116.8333
$ws.Range("L2").Value = 3863.6667
$ws.Range("M2").Value = -3003.8333
$ws.Range("N2").Value = -4089.6667
$ws.Range("H32").Value = 3163.7656
$ws.Range("I32").Value = 3104.5322
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 3104.5322
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -2817.5322
$ws.Range("N32").Value = -5574
$ws.Range("H61").Value = 7206713.5
$ws.Range("I61").Value = 9155926
$ws.Range("J61").Value = 2225393.5
$ws.Range("K61").Value = 9155926
$ws.Range("L61").Value = 2225393.5
$ws.Range("M61").Value = -9155714
$ws.Range("N61").Value = -2225817.5
$ws.Range("H116").Value = 3303.5417
$ws.Range("I116").Value = 3116.8333
$ws.Range("J116").Value = 3863.6667
$ws.Range("K116").Value = 3116.8333
$ws.Range("L116").Value = 3863.6667
$ws.Range("M116").Value = -822.8332999999998
$ws.Range("N116").Value = -8451.6667
$ws.Range("H122").Value = 3534.7908
$ws.Range("I122").Value = 3419.389
$ws.Range("J122").Value = 4128.2856
$ws.Range("K122").Value = 10258.167
$ws.Range("L122").Value = 12384.8568
$ws.Range("M122").Value = -7808.167000000001
$ws.Range("N122").Value = -17284.8568
$ws.Range("H132").Value = 1855743.6
$ws.Range("I132").Value = 3646.8635
$ws.Range("J132").Value = 10004970
$ws.Range("K132").Value = 10940.5905
$ws.Range("L132").Value = 30014910
$ws.Range("M132").Value = -8410.5905
$ws.Range("N132").Value = -30019970
$ws.Range("H133").Value = 64444
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 64444
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 64444
$ws.Range("N133").Value = -69504
$ws.Range("H136").Value = 7206713.5
$ws.Range("I136").Value = 9155926
$ws.Range("J136").Value = 2225393.5
$ws.Range("K136").Value = 27467778
$ws.Range("L136").Value = 6676180.5
$ws.Range("M136").Value = -27465228
$ws.Range("N136").Value = -6681280.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3303.5417
$ws.Range("I3").Value = 3116.8333
$ws.Range("J3").Value = 3863.6667
$ws.Range("K3").Value = 3116.8333
$ws.Range("L3").Value = 3863.6667
$ws.Range("M3").Value = -3002.8333
$ws.Range("N3").Value = -4091.6667
$ws.Range("H107").Value = 2684.9412
$ws.Range("I107").Value = 3038.4
$ws.Range("J107").Value = 1703.1111
$ws.Range("K107").Value = 3038.4
$ws.Range("L107").Value = 1703.1111
$ws.Range("M107").Value = -1118.4
$ws.Range("N107").Value = -5543.1111
$ws.Range("H132").Value = 192500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 192500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 192500
$ws.Range("N132").Value = -202620
$ws.Range("H134").Value = 7694495.5
$ws.Range("I134").Value = 2312.9092
$ws.Range("J134").Value = 50001500
$ws.Range("K134").Value = 6938.7276
$ws.Range("L134").Value = 150004500
$ws.Range("M134").Value = -4403.7276
$ws.Range("N134").Value = -150009570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16952532
$ws.Range("I31").Value = 30306378
$ws.Range("J31").Value = 3419.1155
$ws.Range("K31").Value = 30306378
$ws.Range("L31").Value = 3419.1155
$ws.Range("M31").Value = -30306083
$ws.Range("N31").Value = -4009.1155
$ws.Range("H34").Value = 16952532
$ws.Range("I34").Value = 30306378
$ws.Range("J34").Value = 3419.1155
$ws.Range("K34").Value = 30306378
$ws.Range("L34").Value = 3419.1155
$ws.Range("M34").Value = -30306176
$ws.Range("N34").Value = -3823.1155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1508154.1
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1508154.1
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 4524462.300000001
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4525028.300000001
$ws.Range("H34").Value = 5817.8335
$ws.Range("I34").Value = 266
$ws.Range("J34").Value = 11369.667
$ws.Range("K34").Value = 798
$ws.Range("L34").Value = 34109.001
$ws.Range("M34").Value = -714
$ws.Range("N34").Value = -34277.001
$ws.Range("H39").Value = 33333
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 33333
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 99999
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -100587
$ws.Range("H55").Value = 4229.7856
$ws.Range("I55").Value = 2034.909
$ws.Range("J55").Value = 12277.667
$ws.Range("K55").Value = 6104.727000000001
$ws.Range("L55").Value = 36833.001
$ws.Range("M55").Value = -5927.727000000001
$ws.Range("N55").Value = -37187.001
$ws.Range("H68").Value = 3029.7693
$ws.Range("I68").Value = 2732.6667
$ws.Range("J68").Value = 3284.4285
$ws.Range("K68").Value = 8198.000100000001
$ws.Range("L68").Value = 9853.2855
$ws.Range("M68").Value = -7387.000100000001
$ws.Range("N68").Value = -11475.2855
$ws.Range("H71").Value = 3029.7693
$ws.Range("I71").Value = 2732.6667
$ws.Range("J71").Value = 3284.4285
$ws.Range("K71").Value = 24594.0003
$ws.Range("L71").Value = 29559.8565
$ws.Range("M71").Value = -20538.0003
$ws.Range("N71").Value = -37671.8565
$ws.Range("H129").Value = 3945.7368
$ws.Range("I129").Value = 2655.4167
$ws.Range("J129").Value = 6157.7144
$ws.Range("K129").Value = 7966.250100000001
$ws.Range("L129").Value = 18473.1432
$ws.Range("M129").Value = -2966.250100000001
$ws.Range("N129").Value = -28473.1432
$ws.Range("H131").Value = 3422.348
$ws.Range("I131").Value = 1991.75
$ws.Range("J131").Value = 4983
$ws.Range("K131").Value = 5975.25
$ws.Range("L131").Value = 14949
$ws.Range("M131").Value = -935.25
$ws.Range("N131").Value = -25029

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2550.7222
$ws.Range("I80").Value = 1732.8334
$ws.Range("J80").Value = 4186.5
$ws.Range("K80").Value = 1732.8334
$ws.Range("L80").Value = 4186.5
$ws.Range("M80").Value = -734.8334
$ws.Range("N80").Value = -6182.5
$ws.Range("H83").Value = 2550.7222
$ws.Range("I83").Value = 1732.8334
$ws.Range("J83").Value = 4186.5
$ws.Range("K83").Value = 8664.166999999999
$ws.Range("L83").Value = 20932.5
$ws.Range("M83").Value = -3672.166999999999
$ws.Range("N83").Value = -30916.5
$ws.Range("H102").Value = 2260.2104
$ws.Range("I102").Value = 1896
$ws.Range("J102").Value = 3049.3333
$ws.Range("K102").Value = 1896
$ws.Range("L102").Value = 3049.3333
$ws.Range("M102").Value = -274
$ws.Range("N102").Value = -6293.3333
$ws.Range("H126").Value = 1720.4445
$ws.Range("I126").Value = 1723
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 5169
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -2699
$ws.Range("N126").Value = -10040
$ws.Range("H132").Value = 5026479.5
$ws.Range("I132").Value = 2497.0715
$ws.Range("J132").Value = 19093632
$ws.Range("K132").Value = 7491.2145
$ws.Range("L132").Value = 57280896
$ws.Range("M132").Value = -4961.2145
$ws.Range("N132").Value = -57285956
$ws.Range("H141").Value = 78404.836
$ws.Range("I141").Value = 67000
$ws.Range("J141").Value = 80685.8
$ws.Range("K141").Value = 67000
$ws.Range("L141").Value = 80685.8
$ws.Range("M141").Value = -61820
$ws.Range("N141").Value = -91045.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5613.2
$ws.Range("I7").Value = 4903.6665
$ws.Range("J7").Value = 11999
$ws.Range("K7").Value = 4903.6665
$ws.Range("L7").Value = 11999
$ws.Range("M7").Value = -4791.6665
$ws.Range("N7").Value = -12223
$ws.Range("H22").Value = 18862842
$ws.Range("I22").Value = 18862842
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 18862842
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -18862547
$ws.Range("H27").Value = 18862842
$ws.Range("I27").Value = 18862842
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 18862842
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -18862735
$ws.Range("H40").Value = 2629.2856
$ws.Range("I40").Value = 2650.8
$ws.Range("J40").Value = 2199
$ws.Range("K40").Value = 2650.8
$ws.Range("L40").Value = 2199
$ws.Range("M40").Value = -2514.8
$ws.Range("N40").Value = -2471
$ws.Range("H93").Value = 2928875.5
$ws.Range("I93").Value = 1908.7
$ws.Range("J93").Value = 6181061
$ws.Range("K93").Value = 1908.7
$ws.Range("L93").Value = 6181061
$ws.Range("M93").Value = -660.7
$ws.Range("N93").Value = -6183557
$ws.Range("H126").Value = 5613.2
$ws.Range("I126").Value = 4903.6665
$ws.Range("J126").Value = 11999
$ws.Range("K126").Value = 14710.9995
$ws.Range("L126").Value = 35997
$ws.Range("M126").Value = -12240.9995
$ws.Range("N126").Value = -40937
$ws.Range("H132").Value = 2885.5386
$ws.Range("I132").Value = 1578
$ws.Range("J132").Value = 6434.5713
$ws.Range("K132").Value = 4734
$ws.Range("L132").Value = 19303.7139
$ws.Range("M132").Value = -2204
$ws.Range("N132").Value = -24363.7139
$ws.Range("H136").Value = 3538.8333
$ws.Range("I136").Value = 3538.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10616.4999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8066.499899999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4509.8
$ws.Range("I122").Value = 4183.3335
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 12550.0005
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -10100.0005
$ws.Range("N122").Value = -19898.5
$ws.Range("H126").Value = 4413.3
$ws.Range("I126").Value = 4966.625
$ws.Range("J126").Value = 2200
$ws.Range("K126").Value = 14899.875
$ws.Range("L126").Value = 6600
$ws.Range("M126").Value = -12429.875
$ws.Range("N126").Value = -11540
$ws.Range("H132").Value = 296360.28
$ws.Range("I132").Value = 2394.6538
$ws.Range("J132").Value = 1251748.6
$ws.Range("K132").Value = 7183.9614
$ws.Range("L132").Value = 3755245.8
$ws.Range("M132").Value = -4653.9614
$ws.Range("N132").Value = -3760305.8
$ws.Range("H141").Value = 177129.8
$ws.Range("I141").Value = 87650
$ws.Range("J141").Value = 199499.75
$ws.Range("K141").Value = 87650
$ws.Range("L141").Value = 199499.75
$ws.Range("M141").Value = -82470
$ws.Range("N141").Value = -209859.75
